$wb = $excel.ActiveWorkbook

# "Repayment schedule" is the 3rd sheet (Input, Summary, Repayment schedule, Transactions)
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the old "Late" column (column N), shifting
# Late / heading(Original) / Outstanding one column to the right.
$wsRepay.Columns("N:N").Insert()

# The inserted column picks up the column-M (left neighbour) width in real
# Excel; reproduce that explicit width on the new column.
$wsRepay.Columns("N:N").ColumnWidth = 9.8

# Make "Repayment schedule" the active sheet/tab, with L14 selected - this
# also clears the previously active selection/tab flag on "Transactions".
$wsRepay.Activate()
$wsRepay.Range("L14").Select()
